$d = $word.ActiveDocument

$rHeading = $d.Content
$rHeading.Find.Execute("Delete Report")
$headingStart = $rHeading.Start
$rTbl = $d.Tables.Item(4)
$tblStart = $rTbl.Range.Start
$delRange = $d.Range($headingStart, $tblStart)
$delRange.Delete()
$d.Tables.Item(4).Delete()

$bmRange = $d.Range(2559, 2559)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"></pkg:package>'
$bmRange.InsertXML("<w:bookmarkStart xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:id='0' w:name='_GoBack'/><w:bookmarkEnd xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:id='0'/>")
